$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new test rows for the downside_frequency / upside_frequency macros,
# appended right after the existing last row (118).
$ws.Range("A119").Value = "downside frequency"
$ws.Range("B119").Value = "Test downside frequency"
$ws.Range("C119").Value = "downside_frequency_test"

$ws.Range("A120").Value = "upside frequency"
$ws.Range("B120").Value = "Test upside frequency"
$ws.Range("C120").Value = "upside_frequency_test"

$ws.Range("C120").Select()

